$wb = $excel.ActiveWorkbook

# --- Permissions sheet: rename table target from "Permissions" to "Permission" ---
$wsPermission = $wb.Worksheets.Item("Permissions")
$wsPermission.Range("A1").Value = "Permission"
$wsPermission.Range("E3:E17").Select()

# --- Genders sheet: rename table target from "Genders" to "Gender" and drop the
#     NameLatin column (column C), shifting the formula column from D to C ---
$wsGender = $wb.Worksheets.Item("Genders")
$wsGender.Range("A1").Value = "Gender"
$wsGender.Columns.Item(3).Delete()
$wsGender.Range("C3").Formula = '=CONCATENATE("insert into ",$A$1,"(",$B$2,") values(N''",B2,"'');")'
$wsGender.Range("C4").Formula = '=CONCATENATE("insert into ",$A$1,"(",$B$2,") values(N''",B3,"'');")'

# --- PartnerPermission sheet: no data changes, just move the selection ---
$wsPartnerPermission = $wb.Worksheets.Item("PartnerPermission")
$wsPartnerPermission.Range("E3:E27").Select()

# --- Make Genders the active sheet/tab (was PartnerPermission) ---
$wsGender.Activate()
$wsGender.Range("D25").Select()
